# Updated cryptos list with refreshed Price (D) and Volume(1h) (E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.436.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +1.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.013"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4764"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3709"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07475"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8874"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.848.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07359"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.489"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.598"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.016"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008857"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.014"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.435.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.348"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.069.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.908"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.172"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.281"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7590"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.575"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.954"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.108"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05364"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01966"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.004"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.313"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5364"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.377"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1668"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.561"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4983"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.687"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06326"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
